# Populate the bug-tracker grid (Title / Steps / Attachments / Device / Network / Severity / Priority)
# for bugs b1..b6, add the Notes row, extend the sheet with a couple of blank rows, and adjust
# row heights / column width / selection to match the authored worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (b1)
$ws.Range('B2').Value = 'App crashes when opening the notifications from mobile notifications bar while App not opened'
$ws.Range('C2').Value = '1- Download Yallakora for Play store.                                                        2- Launch the App then close it wait a while you will receive notifications in the phone notifications bar.                                    3- Click the received notification. I didn''t open.                                                    '
$ws.Range('D2').Value = 'Check the attached file called b1'
$ws.Range('E2').Value = 'Samsung Galaxy j4 with Android 9'
$ws.Range('F2').Value = 'WIFI'
$ws.Range('G2').Value = 'Critical'
$ws.Range('H2').Value = 'High'
$ws.Rows.Item(2).RowHeight = 120

# Row 3 (b2)
$ws.Range('B3').Value = 'in Your teams tab Add your team button disappear in portrait mode'
$ws.Range('C3').Value = '1- Launch Yallakora App.                        2- Click your teams tab.                     3- Change your phone to Portrait mode.                                                       Add your team button disappeared.'
$ws.Range('D3').Value = 'Check the attached file called b2'
$ws.Range('E3').Value = 'Samsung Galaxy j4 with Android 9'
$ws.Range('F3').Value = 'WIFI'
$ws.Range('G3').Value = 'High'
$ws.Range('H3').Value = 'High'
$ws.Rows.Item(3).RowHeight = 90

# Row 4 (b3)
$ws.Range('B4').Value = 'Any team name which contains two words in vertical mode the second word is missing'
$ws.Range('C4').Value = '1- Launch Yallakora App.                        2- Click matches tab.                                      3- Scroll down to the European matches.                                                    4- Click north irland match.             Word north is missing and when clicking on the match finished lab not fully appeared and there is no data in formation.'
$ws.Range('D4').Value = 'Check the attached file called b3'
$ws.Range('E4').Value = 'Samsung Galaxy j4 with Android 9'
$ws.Range('F4').Value = 'WIFI'
$ws.Range('G4').Value = 'Low'
$ws.Range('H4').Value = 'Medium '
$ws.Rows.Item(4).RowHeight = 135

# Row 5 (b4)
$ws.Range('B5').Value = 'There is two results for a finished match between Cheli and Gana'
$ws.Range('C5').Value = '1- Launch Yallakora App.                        2- Click matches tab.                             3- Select Tuesday 14 June, First finished match between cheli and Gana there is two results (0-0 and 3-1)    '
$ws.Range('D5').Value = 'Check the attached file called b4'
$ws.Range('E5').Value = 'Samsung Galaxy j4 with Android 9'
$ws.Range('F5').Value = 'WIFI'
$ws.Range('G5').Value = 'Medium as when you click the match it shows the right result'
$ws.Range('H5').Value = 'High as this is a kora App news'
$ws.Rows.Item(5).RowHeight = 90

# Row 6 (b5)
$ws.Range('B6').Value = 'The update teams sign is active also I didn''t add teams yet.'
$ws.Range('C6').Value = '1- Launch Yallakora App.                        2- Click Your teams tab.                     Also I didn''t add teams but the upper right side update sign is active.  '
$ws.Range('D6').Value = 'Check the attached file called b5'
$ws.Range('E6').Value = 'Samsung Galaxy j4 with Android 9'
$ws.Range('F6').Value = 'WIFI'
$ws.Range('G6').Value = 'Low as when click the update sign it opens to selcet team '
$ws.Range('H6').Value = 'Low as the update sign color the same if it active or inactive'
$ws.Rows.Item(6).RowHeight = 90

# Row 7 (b6)
$ws.Range('B7').Value = 'There is invalid date word for first item of important news in home (as it an Ad not a news) and sharing it  take about two seconds to open sharing options'
$ws.Range('C7').Value = '1- Launch Yallakora App.                        2- Click Home tab.                                3- Scroll down to Important news check first item (An Ad)it''s has invaild date and when clicking it and click the upper right side sharing sign it takes about two seconds (this more than usual )to open sharing options.'
$ws.Range('D7').Value = 'Check the attached file called b6'
$ws.Range('E7').Value = 'Samsung Galaxy j4 with Android 9'
$ws.Range('F7').Value = 'WIFI'
$ws.Range('G7').Value = 'Critical as this is the Home tab'
$ws.Range('H7').Value = 'High'
$ws.Rows.Item(7).RowHeight = 150

# Row 11: Notes label + bold, wrapped note text; clear the rest of the row
$ws.Range('A11').Value = 'Notes:'
$ws.Range('B11').Value = ' App works correctly in cases (receiving a phone call, SMS, low battery, high screen brightness, and low screen brightness)'
$ws.Range('B11').Font.Bold = $true
$ws.Range('B11').WrapText = $true
$ws.Range('C11:H11').Clear()
$ws.Rows.Item(11).RowHeight = 90

# Extend the used range with two more blank, formatted rows (14-15)
$ws.Range('A14:H15').WrapText = $true

# Widen column H slightly to fit the new Priority values
$ws.Columns.Item(8).ColumnWidth = 9.86

# Select the header row, as in the authored workbook
[void]$ws.Rows.Item(1).Select()
